# "fix lại báo cáo export"
# Replace the "sau thuế" (after-tax) columns with "sau giảm" (after-promotion)
# columns in the Direct Sales Order General report template.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 7): labels for the three summary columns J/K/L ---
# J7: "Tổng tiền sau thuế" -> "Doanh thu"
# K7: "Giảm giá" (unchanged)
# L7: "Tôgnr tiền" (typo) -> "Tổng tiền sau giảm"
$ws.Range("J7").Value = "Doanh thu"
$ws.Range("K7").Value = "Giảm giá"
$ws.Range("L7").Value = "Tổng tiền sau giảm"

# --- Placeholder row (row 9): per-order merge tags for J/K/L ---
# J9: {{...TotalAfterTax}} -> {{...Total}}
# K9: {{...TotalAfterTax}}-ish slot now holds PromotionValue
# L9: {{...Total}} -> {{...TotalAfterPromotion}}
$ws.Range("J9").Value = "{{ReportSalesOrderGenerals.SalesOrders.Total}}"
$ws.Range("K9").Value = "{{ReportSalesOrderGenerals.SalesOrders.PromotionValue}}"
$ws.Range("L9").Value = "{{ReportSalesOrderGenerals.SalesOrders.TotalAfterPromotion}}"

# --- Totals row (row 10): grand-total merge tags for J/K/L ---
# J10: {{Total.TotalAfterTax}} -> {{Total.Total}}
# K10: {{Total.TotalAfterTax}}-ish slot now holds {{Total.PromotionValue}}
# L10: {{Total.Total}} -> {{Total.TotalAfterPromotion}}
$ws.Range("J10").Value = "{{Total.Total}}"
$ws.Range("K10").Value = "{{Total.PromotionValue}}"
$ws.Range("L10").Value = "{{Total.TotalAfterPromotion}}"

# --- Selection moved (cosmetic, reflects where the author left the cursor) ---
$ws.Range("F12").Select()
